function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue 2 4 '19.928.89'
Set-TextValue 2 5 '  -8.30%  '

Set-TextValue 3 4 '1.401.85'
Set-TextValue 3 5 '  -8.98%  '

Set-TextValue 4 4 '1.003'
Set-TextValue 4 5 '  +0.22%  '

Set-TextValue 5 5 '  +0.23%  '

Set-TextValue 6 4 '272.55'
Set-TextValue 6 5 '  -6.06%  '

Set-TextValue 7 4 '0.3704'
Set-TextValue 7 5 '  -4.86%  '

Set-TextValue 8 4 '0.3063'
Set-TextValue 8 5 '  -3.79%  '

Set-TextValue 9 4 '39.07'
Set-TextValue 9 5 '  -9.37%  '

Set-TextValue 10 4 '0.9912'
Set-TextValue 10 5 '  -6.46%  '

Set-TextValue 11 4 '0.06544'
Set-TextValue 11 5 '  -9.04%  '

Set-TextValue 12 5 '  +0.33%  '

Set-TextValue 13 4 '5.402'
Set-TextValue 13 5 '  -4.22%  '

Set-TextValue 14 4 '6.147'
Set-TextValue 14 5 '  -7.28%  '

Set-TextValue 15 4 '16.80'
Set-TextValue 15 5 '  -9.94%  '

Set-TextValue 16 4 '1.407.41'
Set-TextValue 16 5 '  -8.67%  '

Set-TextValue 17 5 '  -9.06%  '

Set-TextValue 18 4 '0.05751'
Set-TextValue 18 5 '  -12.58%  '

Set-TextValue 19 4 '73.36'
Set-TextValue 19 5 '  -11.81%  '

Set-TextValue 20 5 '  +0.29%  '

Set-TextValue 21 4 '5.571'
Set-TextValue 21 5 '  -9.42%  '

Set-TextValue 22 4 '14.35'
Set-TextValue 22 5 '  -6.70%  '

Set-TextValue 23 4 '10.74'
Set-TextValue 23 5 '  -1.21%  '

Set-TextValue 24 4 '2.327'
Set-TextValue 24 5 '  -3.10%  '

Set-TextValue 25 4 '19.936.61'
Set-TextValue 25 5 '  -8.25%  '

Set-TextValue 26 5 '  -4.95%  '

Set-TextValue 27 4 '138.29'
Set-TextValue 27 5 '  -5.57%  '

Set-TextValue 28 5 '  -8.19%  '

Set-TextValue 29 4 '1.570.36'
Set-TextValue 29 5 '  -8.44%  '

Set-TextValue 30 4 '108.64'
Set-TextValue 30 5 '  -7.51%  '

Set-TextValue 31 4 '3.825'
Set-TextValue 31 5 '  -21.02%  '

Set-TextValue 32 4 '5.362'
Set-TextValue 32 5 '  -9.35%  '

Set-TextValue 33 4 '0.8411'
Set-TextValue 33 5 '  -13.00%  '

Set-TextValue 34 4 '0.07700'
Set-TextValue 34 5 '  -6.13%  '

Set-TextValue 35 4 '8.396'
Set-TextValue 35 5 '  -4.80%  '

Set-TextValue 36 2 'InternetComputer(DFINITY)'
Set-TextValue 36 3 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 36 4 '4.798'
Set-TextValue 36 5 '  -6.39%  '

Set-TextValue 37 2 'Hedera'
Set-TextValue 37 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 37 4 '0.05738'
Set-TextValue 37 5 '  -5.79%  '

Set-TextValue 38 2 'Frax'
Set-TextValue 38 3 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 38 4 '1.003'
Set-TextValue 38 5 '  +0.30%  '

Set-TextValue 39 4 '0.1908'
Set-TextValue 39 5 '  -6.42%  '

Set-TextValue 40 4 '0.02024'
Set-TextValue 40 5 '  -7.94%  '

Set-TextValue 41 4 '10.19'
Set-TextValue 41 5 '  -4.36%  '

Set-TextValue 42 5 '  -11.24%  '

Set-TextValue 43 4 '1.266'
Set-TextValue 43 5 '  -12.66%  '

Set-TextValue 44 4 '0.5264'
Set-TextValue 44 5 '  -8.16%  '

Set-TextValue 45 4 '3.523'
Set-TextValue 45 5 '  -5.81%  '

Set-TextValue 46 4 '12.14'
Set-TextValue 46 5 '  -7.40%  '

Set-TextValue 47 4 '0.5085'
Set-TextValue 47 5 '  -7.33%  '

Set-TextValue 48 4 '1.798'
Set-TextValue 48 5 '  -4.00%  '

Set-TextValue 49 4 '108.95'
Set-TextValue 49 5 '  -6.37%  '

Set-TextValue 50 5 '  -9.89%  '

Set-TextValue 51 4 '1.004'
Set-TextValue 51 5 '  +0.32%  '
